$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.416.64"
$ws.Range("E2").Value = "  -0.07%  "

$ws.Range("D3").Value = "3.144.53"
$ws.Range("E3").Value = "  -0.35%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "611.13"
$ws.Range("E5").Value = "  +0.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.60"
$ws.Range("E6").Value = "  -1.89%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "3.141.92"
$ws.Range("E8").Value = "  -0.25%  "

$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("E10").Value = "  -0.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.36"
$ws.Range("E11").Value = "  -3.42%  "

$ws.Range("E12").Value = "  -0.34%  "

$ws.Range("E13").Value = "  +0.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.50"
$ws.Range("E14").Value = "  -1.40%  "

$ws.Range("D15").Value = "3.664.12"
$ws.Range("E15").Value = "  -0.29%  "

$ws.Range("E16").Value = "  +2.99%  "

$ws.Range("D17").Value = "64.375.25"
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("D18").Value = "3.143.49"
$ws.Range("E18").Value = "  -0.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.84"
$ws.Range("E19").Value = "  -1.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "476.44"
$ws.Range("E20").Value = "  -0.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.65"
$ws.Range("E21").Value = "  +0.47%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.721"
$ws.Range("E22").Value = "  +1.59%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.82"
$ws.Range("E23").Value = "  +1.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.63"
$ws.Range("E24").Value = "  -0.89%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.84"
$ws.Range("E25").Value = "  +1.70%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("E27").Value = "  -3.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.58"
$ws.Range("E28").Value = "  +2.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.41"
$ws.Range("E29").Value = "  +9.07%  "

$ws.Range("E30").Value = "  +3.05%  "

$ws.Range("E31").Value = "  -4.57%  "

$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.61"
$ws.Range("E33").Value = "  +1.73%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.64"
$ws.Range("E34").Value = "  -3.92%  "

$ws.Range("E35").Value = "  +0.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.93"
$ws.Range("E36").Value = "  -1.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.58"
$ws.Range("E37").Value = "  -2.93%  "

$ws.Range("D38").Value = "0.0₃0742"
$ws.Range("E38").Value = "  +3.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.01"
$ws.Range("E39").Value = "  +3.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "451.63"
$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0394"
$ws.Range("E41").Value = "  -0.59%  "

$ws.Range("E42").Value = "  -0.24%  "

$ws.Range("E43").Value = "  -1.26%  "

$ws.Range("D44").Value = "2.848.49"
$ws.Range("E44").Value = "  +0.31%  "

$ws.Range("E45").Value = "  -0.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.26"
$ws.Range("E46").Value = "  +0.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.42"
$ws.Range("E47").Value = "  +4.70%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.41"
$ws.Range("E48").Value = "  +0.05%  "

$ws.Range("E49").Value = "  +0.10%  "

$ws.Range("E50").Value = "  +0.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.27"
$ws.Range("E51").Value = "  +1.57%  "
